$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 inherits row 2's formatting (incl. the Text number format on the
# ID/phone columns) before we fill in its values, so numeric-looking
# strings stay text instead of being auto-converted to numbers.
$ws.Range("A2:L2").Copy()
$ws.Range("A3:L3").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(2).RowHeight

# Row 3 = second unit of the same order as (the original) row 2.
$ws.Range("A3").Value = "23020076246"
$ws.Range("B3").Value = "199201934887701"
$ws.Range("C3").Value = 465
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "2013-01-07 10:52"
$ws.Range("G3").Value = "江苏苏州市昆山市千灯镇少卿西路卿峰丽景南区苏式建筑198幢"
$ws.Range("I3").Value = "益菱"
$ws.Range("J3").Value = "13862393981"
$ws.Range("L3").Value = "江苏苏州市昆山市(周庄、锦溪、淀山湖、千灯、张浦、陆家、花桥)江苏苏州市昆山市(周庄、锦溪、淀山湖、千灯、张浦、陆家、花桥)周庄锦溪淀山湖千灯张浦陆家花桥千灯镇少卿西路卿峰丽景南区苏式建筑198幢"

# Row 2 keeps the rest of its data but gets its own distinct order number.
$ws.Range("A2").Value = "23020076245"

$ws.PageSetup.PaperSize = 9 # xlPaperA4
$ws.PageSetup.Orientation = 1 # xlPortrait

$ws.Range("A2").Select()
